$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue "D2" "29.951.56"
Set-TextValue "E2" "  +0.32%  "
Set-TextValue "D3" "1.891.54"
Set-TextValue "E3" "  -0.30%  "
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "0.8213"
Set-TextValue "E5" "  +5.19%  "
Set-TextValue "D6" "241.42"
Set-TextValue "E6" "  +0.54%  "
Set-TextValue "E7" "  +0.08%  "
Set-TextValue "E8" "  +5.57%  "
Set-TextValue "D9" "26.40"
Set-TextValue "E9" "  +3.34%  "
Set-TextValue "D10" "0.07023"
Set-TextValue "E10" "  +2.31%  "
Set-TextValue "D11" "0.08040"
Set-TextValue "E11" "  +0.70%  "
Set-TextValue "D12" "0.7457"
Set-TextValue "E12" "  +1.01%  "
Set-TextValue "D13" "1.889.57"
Set-TextValue "E13" "  -1.55%  "
Set-TextValue "D14" "5.203"
Set-TextValue "E14" "  +0.45%  "
Set-TextValue "D15" "92.23"
Set-TextValue "E15" "  +0.81%  "
Set-TextValue "D16" "29.946.86"
Set-TextValue "E16" "  +0.28%  "
Set-TextValue "D17" "14.03"
Set-TextValue "E17" "  +1.63%  "
Set-TextValue "D18" "5.885"
Set-TextValue "E18" "  -0.24%  "
Set-TextValue "D19" "244.63"
Set-TextValue "E19" "  +0.16%  "
Set-TextValue "D20" "0.000007760"
Set-TextValue "E20" "  +0.63%  "
Set-TextValue "E21" "  +0.07%  "
Set-TextValue "D22" "2.141.49"
Set-TextValue "E22" "  -1.28%  "
Set-TextValue "E23" "  +0.04%  "
Set-TextValue "D24" "6.927"
Set-TextValue "E24" "  -0.04%  "
Set-TextValue "D25" "0.1558"
Set-TextValue "E25" "  +16.93%  "
Set-TextValue "D26" "166.06"
Set-TextValue "E26" "  -0.48%  "
Set-TextValue "D27" "9.196"
Set-TextValue "E27" "  -0.99%  "
Set-TextValue "D28" "18.82"
Set-TextValue "D29" "2.084"
Set-TextValue "E29" "  +2.76%  "
Set-TextValue "E30" "  -1.68%  "
Set-TextValue "D31" "1.516"
Set-TextValue "E31" "  +0.41%  "
Set-TextValue "E32" "  -0.19%  "
Set-TextValue "D33" "0.05637"
Set-TextValue "E33" "  +6.94%  "
Set-TextValue "D34" "4.066"
Set-TextValue "E34" "  -0.17%  "
Set-TextValue "E35" "  +2.10%  "
Set-TextValue "D36" "0.7271"
Set-TextValue "E36" "  -0.27%  "
Set-TextValue "E37" "  -0.39%  "
Set-TextValue "D38" "0.01912"
Set-TextValue "E38" "  +0.30%  "
Set-TextValue "D39" "2.779"
Set-TextValue "E39" "  +0.20%  "
Set-TextValue "D40" "0.4424"
Set-TextValue "E40" "  +0.04%  "
Set-TextValue "D41" "71.84"
Set-TextValue "E41" "  -0.33%  "
Set-TextValue "D42" "5.954"
Set-TextValue "E42" "  -3.84%  "
Set-TextValue "D43" "0.8437"
Set-TextValue "E43" "  +1.06%  "
Set-TextValue "D44" "1.001"
Set-TextValue "E44" "  +0.04%  "
Set-TextValue "D45" "1.870"
Set-TextValue "E45" "  -0.44%  "
Set-TextValue "B46" "Quant"
Set-TextValue "C46" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D46" "100.58"
Set-TextValue "E46" "  +0.18%  "
Set-TextValue "B47" "Aptos"
Set-TextValue "C47" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D47" "7.564"
Set-TextValue "E47" "  -0.41%  "
Set-TextValue "D48" "9.713"
Set-TextValue "E48" "  -0.65%  "
Set-TextValue "D49" "991.00"
Set-TextValue "E49" "  +7.28%  "
Set-TextValue "D50" "2.039.34"
Set-TextValue "E50" "  -1.06%  "
Set-TextValue "D51" "35.92"
Set-TextValue "E51" "  -0.42%  "
